{"js": "// This table holds simple arithmetic problems (e.g. \"53-27=26\"). The edit\n// replaces each cell's text with a new problem/answer string, preserving\n// cell order, row/column shape, and all existing formatting (fonts, size,\n// alignment, etc. are untouched since only the text content changes).\n// Some \"old\" strings repeat (e.g. \"31+62=93\" appears twice but maps to two\n// different replacements), so cells are matched strictly by position/order\n// rather than by a global find-and-replace of the old text.\nconst pairs = [\n  [\"53-27=26\", \"55-45=10\"],\n  [\"0+10=10\", \"29+50=79\"],\n  [\"29-5=24\", \"23+37=60\"],\n  [\"96-67=29\", \"77-18=59\"],\n  [\"98-67=31\", \"23+57=80\"],\n  [\"31+62=93\", \"9+7=16\"],\n  [\"17+82=99\", \"25+12=37\"],\n  [\"15+83=98\", \"20+13=33\"],\n  [\"56-53=3\", \"31-13=18\"],\n  [\"12+83=95\", \"10+82=92\"],\n  [\"11+47=58\", \"81-18=63\"],\n  [\"83-25=58\", \"33-15=18\"],\n  [\"76-48=28\", \"12+50=62\"],\n  [\"70+13=83\", \"20-16=4\"],\n  [\"90-60=30\", \"13+61=74\"],\n  [\"2+51=53\", \"28+32=60\"],\n  [\"38-8=30\", \"91-49=42\"],\n  [\"40+11=51\", \"43+23=66\"],\n  [\"35+2=37\", \"0+88=88\"],\n  [\"25+62=87\", \"25+20=45\"],\n  [\"81-42=39\", \"39-1=38\"],\n  [\"37+45=82\", \"9+86=95\"],\n  [\"14+23=37\", \"16+18=34\"],\n  [\"7+79=86\", \"70-15=55\"],\n  [\"67+30=97\", \"79-57=22\"],\n  [\"71-66=5\", \"34-18=16\"],\n  [\"3+48=51\", \"85+2=87\"],\n  [\"30+63=93\", \"7+42=49\"],\n  [\"53-41=12\", \"72+26=98\"],\n  [\"86-65=21\", \"35+3=38\"],\n  [\"87-38=49\", \"56-32=24\"],\n  [\"92-82=10\", \"50+39=89\"],\n  [\"22+77=99\", \"33-31=2\"],\n  [\"75-15=60\", \"85-83=2\"],\n  [\"6+92=98\", \"75-34=41\"],\n  [\"81-22=59\", \"98-71=27\"],\n  [\"95-74=21\", \"51-26=25\"],\n  [\"94-70=24\", \"51-32=19\"],\n  [\"61-51=10\", \"8+34=42\"],\n  [\"84-54=30\", \"59+32=91\"],\n  [\"40+21=61\", \"87-58=29\"],\n  [\"37+17=54\", \"46-26=20\"],\n  [\"41-24=17\", \"69+7=76\"],\n  [\"25+33=58\", \"49+13=62\"],\n  [\"39+43=82\", \"78-28=50\"],\n  [\"14+16=30\", \"40-5=35\"],\n  [\"31+62=93\", \"12+41=53\"],\n  [\"89-10=79\", \"60-42=18\"],\n  [\"86-54=32\", \"65+30=95\"],\n  [\"81-75=6\", \"59+13=72\"],\n  [\"88-80=8\", \"76+17=93\"],\n  [\"36+60=96\", \"85-24=61\"],\n  [\"0+94=94\", \"87-83=4\"],\n  [\"59+11=70\", \"85-56=29\"],\n  [\"21+35=56\", \"61-2=59\"],\n  [\"34+43=77\", \"42+29=71\"],\n  [\"2+94=96\", \"85-10=75\"],\n  [\"80+19=99\", \"35-25=10\"],\n  [\"52-44=8\", \"92-56=36\"],\n  [\"40+9=49\", \"9+45=54\"],\n  [\"37+33=70\", \"97-28=69\"],\n  [\"52+31=83\", \"93-17=76\"],\n  [\"67-66=1\", \"58-5=53\"],\n  [\"73+8=81\", \"76-74=2\"],\n  [\"77-16=61\", \"92+4=96\"],\n  [\"49+7=56\", \"68-15=53\"],\n  [\"49-43=6\", \"21+72=93\"],\n  [\"2+70=72\", \"49+38=87\"],\n  [\"59-12=47\", \"34-26=8\"],\n  [\"89-47=42\", \"7+72=79\"],\n  [\"28+62=90\", \"15+70=85\"],\n  [\"41-30=11\", \"60-21=39\"],\n  [\"27+52=79\", \"44+5=49\"],\n  [\"82-60=22\", \"78+12=90\"],\n  [\"87-38=49\", \"52-28=24\"],\n  [\"54+38=92\", \"94-55=39\"],\n  [\"14+1=15\", \"41-20=21\"],\n  [\"64-31=33\", \"51-34=17\"],\n  [\"80-74=6\", \"72+26=98\"],\n  [\"51-13=38\", \"66-46=20\"],\n  [\"29+42=71\", \"29-29=0\"],\n  [\"23-12=11\", \"79-71=8\"],\n  [\"0+71=71\", \"3+21=24\"],\n  [\"42-15=27\", \"68-65=3\"],\n  [\"15+43=58\", \"43-7=36\"],\n  [\"68+29=97\", \"15+38=53\"],\n  [\"90-73=17\", \"38-33=5\"],\n  [\"92-23=69\", \"51-49=2\"],\n  [\"17+28=45\", \"72-39=33\"],\n  [\"53-52=1\", \"16+44=60\"],\n  [\"95-87=8\", \"72-4=68\"],\n  [\"36-19=17\", \"71+21=92\"],\n  [\"15+0=15\", \"79+2=81\"],\n  [\"41-18=23\", \"76+12=88\"],\n  [\"91-87=4\", \"45+9=54\"],\n  [\"61+28=89\", \"9+79=88\"],\n  [\"48+50=98\", \"12+20=32\"],\n  [\"51-23=28\", \"7+68=75\"],\n  [\"74-7=67\", \"29+45=74\"],\n  [\"0+73=73\", \"98-21=77\"]\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in document body\");\n}\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst oldValues = table.values;\nconst flatCount = oldValues.reduce((sum, row) => sum + row.length, 0);\nif (flatCount !== pairs.length) {\n  throw new Error(\n    \"Cell count (\" + flatCount + \") does not match expected pair count (\" + pairs.length + \")\"\n  );\n}\n\n// Walk the grid in row-major order (same order the diff's replacements were\n// recorded in) and build the replacement grid, verifying the pre-edit text\n// along the way.\nconst newValues = [];\nlet idx = 0;\nfor (const row of oldValues) {\n  const newRow = [];\n  for (const cellText of row) {\n    const [expectedOld, replacement] = pairs[idx];\n    if (cellText !== expectedOld) {\n      throw new Error(\n        \"Mismatch at cell \" + idx + \": expected '\" + expectedOld + \"' but found '\" + cellText + \"'\"\n      );\n    }\n    newRow.push(replacement);\n    idx++;\n  }\n  newValues.push(newRow);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# This table holds simple arithmetic problems (e.g. \"53-27=26\"). The edit\n# replaces each cell's text with a new problem/answer string, preserving\n# cell order, row/column shape, and all existing formatting (fonts, size,\n# alignment, etc. are untouched since only the run text changes).\n# Some \"old\" strings repeat (e.g. \"31+62=93\" appears twice but maps to two\n# different replacements), so cells are matched strictly by position/order\n# (row-major, same order as the source diff) rather than by a global\n# find-and-replace of the old text.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$pairs = @(\n  @(\"53-27=26\", \"55-45=10\"),\n  @(\"0+10=10\", \"29+50=79\"),\n  @(\"29-5=24\", \"23+37=60\"),\n  @(\"96-67=29\", \"77-18=59\"),\n  @(\"98-67=31\", \"23+57=80\"),\n  @(\"31+62=93\", \"9+7=16\"),\n  @(\"17+82=99\", \"25+12=37\"),\n  @(\"15+83=98\", \"20+13=33\"),\n  @(\"56-53=3\", \"31-13=18\"),\n  @(\"12+83=95\", \"10+82=92\"),\n  @(\"11+47=58\", \"81-18=63\"),\n  @(\"83-25=58\", \"33-15=18\"),\n  @(\"76-48=28\", \"12+50=62\"),\n  @(\"70+13=83\", \"20-16=4\"),\n  @(\"90-60=30\", \"13+61=74\"),\n  @(\"2+51=53\", \"28+32=60\"),\n  @(\"38-8=30\", \"91-49=42\"),\n  @(\"40+11=51\", \"43+23=66\"),\n  @(\"35+2=37\", \"0+88=88\"),\n  @(\"25+62=87\", \"25+20=45\"),\n  @(\"81-42=39\", \"39-1=38\"),\n  @(\"37+45=82\", \"9+86=95\"),\n  @(\"14+23=37\", \"16+18=34\"),\n  @(\"7+79=86\", \"70-15=55\"),\n  @(\"67+30=97\", \"79-57=22\"),\n  @(\"71-66=5\", \"34-18=16\"),\n  @(\"3+48=51\", \"85+2=87\"),\n  @(\"30+63=93\", \"7+42=49\"),\n  @(\"53-41=12\", \"72+26=98\"),\n  @(\"86-65=21\", \"35+3=38\"),\n  @(\"87-38=49\", \"56-32=24\"),\n  @(\"92-82=10\", \"50+39=89\"),\n  @(\"22+77=99\", \"33-31=2\"),\n  @(\"75-15=60\", \"85-83=2\"),\n  @(\"6+92=98\", \"75-34=41\"),\n  @(\"81-22=59\", \"98-71=27\"),\n  @(\"95-74=21\", \"51-26=25\"),\n  @(\"94-70=24\", \"51-32=19\"),\n  @(\"61-51=10\", \"8+34=42\"),\n  @(\"84-54=30\", \"59+32=91\"),\n  @(\"40+21=61\", \"87-58=29\"),\n  @(\"37+17=54\", \"46-26=20\"),\n  @(\"41-24=17\", \"69+7=76\"),\n  @(\"25+33=58\", \"49+13=62\"),\n  @(\"39+43=82\", \"78-28=50\"),\n  @(\"14+16=30\", \"40-5=35\"),\n  @(\"31+62=93\", \"12+41=53\"),\n  @(\"89-10=79\", \"60-42=18\"),\n  @(\"86-54=32\", \"65+30=95\"),\n  @(\"81-75=6\", \"59+13=72\"),\n  @(\"88-80=8\", \"76+17=93\"),\n  @(\"36+60=96\", \"85-24=61\"),\n  @(\"0+94=94\", \"87-83=4\"),\n  @(\"59+11=70\", \"85-56=29\"),\n  @(\"21+35=56\", \"61-2=59\"),\n  @(\"34+43=77\", \"42+29=71\"),\n  @(\"2+94=96\", \"85-10=75\"),\n  @(\"80+19=99\", \"35-25=10\"),\n  @(\"52-44=8\", \"92-56=36\"),\n  @(\"40+9=49\", \"9+45=54\"),\n  @(\"37+33=70\", \"97-28=69\"),\n  @(\"52+31=83\", \"93-17=76\"),\n  @(\"67-66=1\", \"58-5=53\"),\n  @(\"73+8=81\", \"76-74=2\"),\n  @(\"77-16=61\", \"92+4=96\"),\n  @(\"49+7=56\", \"68-15=53\"),\n  @(\"49-43=6\", \"21+72=93\"),\n  @(\"2+70=72\", \"49+38=87\"),\n  @(\"59-12=47\", \"34-26=8\"),\n  @(\"89-47=42\", \"7+72=79\"),\n  @(\"28+62=90\", \"15+70=85\"),\n  @(\"41-30=11\", \"60-21=39\"),\n  @(\"27+52=79\", \"44+5=49\"),\n  @(\"82-60=22\", \"78+12=90\"),\n  @(\"87-38=49\", \"52-28=24\"),\n  @(\"54+38=92\", \"94-55=39\"),\n  @(\"14+1=15\", \"41-20=21\"),\n  @(\"64-31=33\", \"51-34=17\"),\n  @(\"80-74=6\", \"72+26=98\"),\n  @(\"51-13=38\", \"66-46=20\"),\n  @(\"29+42=71\", \"29-29=0\"),\n  @(\"23-12=11\", \"79-71=8\"),\n  @(\"0+71=71\", \"3+21=24\"),\n  @(\"42-15=27\", \"68-65=3\"),\n  @(\"15+43=58\", \"43-7=36\"),\n  @(\"68+29=97\", \"15+38=53\"),\n  @(\"90-73=17\", \"38-33=5\"),\n  @(\"92-23=69\", \"51-49=2\"),\n  @(\"17+28=45\", \"72-39=33\"),\n  @(\"53-52=1\", \"16+44=60\"),\n  @(\"95-87=8\", \"72-4=68\"),\n  @(\"36-19=17\", \"71+21=92\"),\n  @(\"15+0=15\", \"79+2=81\"),\n  @(\"41-18=23\", \"76+12=88\"),\n  @(\"91-87=4\", \"45+9=54\"),\n  @(\"61+28=89\", \"9+79=88\"),\n  @(\"48+50=98\", \"12+20=32\"),\n  @(\"51-23=28\", \"7+68=75\"),\n  @(\"74-7=67\", \"29+45=74\"),\n  @(\"0+73=73\", \"98-21=77\")\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nif (($rowCount * $colCount) -ne $pairs.Count) {\n    throw \"Cell count ($($rowCount * $colCount)) does not match expected pair count ($($pairs.Count))\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $expectedOld = $pairs[$idx][0]\n        $replacement = $pairs[$idx][1]\n        $cell = $t.Cell($r, $c)\n\n        # TableCell.Range.Text includes trailing cell-mark control characters\n        # (a carriage return followed by the end-of-cell marker); strip them\n        # before comparing against the expected plain text.\n        $current = $cell.Range.Text\n        $current = $current.TrimEnd([char]7, [char]13)\n\n        if ($current -ne $expectedOld) {\n            throw \"Mismatch at row ${r} col ${c}: expected '$expectedOld' but found '$current'\"\n        }\n\n        # Assigning Range.Text replaces only the cell's content and keeps the\n        # end-of-cell marker intact, so the run's formatting (font/size from\n        # rPr) carries over unchanged.\n        $cell.Range.Text = $replacement\n        $idx++\n    }\n}\n"}
